# "INETs in cellular fixed"
# - Region (B2) was wrongly set to APAC, should be EMEA
# - 4G+Cellular (B25) should be False, not True

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")
$ws.Activate()

# Region: APAC -> EMEA
$ws.Range("B2").Value = "EMEA"

# 4G+Cellular: TRUE -> FALSE
$ws.Range("B25").Value = $false

# Update the on-screen selection to match the saved state
$ws.Range("D24").Select()
